$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 362.5
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 425
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 425
$ws.Range("M2").Value = -187
$ws.Range("N2").Value = -651

$ws.Range("H28").Value = 2778.8
$ws.Range("I28").Value = 5447
$ws.Range("K28").Value = 5447
$ws.Range("M28").Value = -4962

$ws.Range("H112").Value = 2095.5
$ws.Range("I112").Value = 640
$ws.Range("J112").Value = 2257.2222
$ws.Range("K112").Value = 1920
$ws.Range("L112").Value = 6771.6666
$ws.Range("M112").Value = -812
$ws.Range("N112").Value = -8987.6666

$ws.Range("H127").Value = 847.75
$ws.Range("I127").Value = 425.9
$ws.Range("K127").Value = 1277.7
$ws.Range("M127").Value = 3682.3

$ws.Range("H132").Value = 190889.39
$ws.Range("I132").Value = 210669.5
$ws.Range("J132").Value = 1000.4
$ws.Range("K132").Value = 632008.5
$ws.Range("L132").Value = 3001.2
$ws.Range("M132").Value = -629478.5
$ws.Range("N132").Value = -8061.2

$ws.Range("H137").Value = 166671660
$ws.Range("I137").Value = 41669416
$ws.Range("J137").Value = 333341340
$ws.Range("K137").Value = 125008248
$ws.Range("L137").Value = 1000024020
$ws.Range("M137").Value = -125005698
$ws.Range("N137").Value = -1000029120

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5548.14
$ws.Range("I32").Value = 5822.125
$ws.Range("J32").Value = 4452.2
$ws.Range("K32").Value = 5822.125
$ws.Range("L32").Value = 4452.2
$ws.Range("M32").Value = -5535.125
$ws.Range("N32").Value = -5026.2

$ws.Range("H45").Value = 1247
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 1247
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 1247
$ws.Range("M45").ClearContents()
$ws.Range("N45").Value = -2001

$ws.Range("H61").Value = 17545672
$ws.Range("I61").Value = 20835136
$ws.Range("J61").Value = 1868
$ws.Range("K61").Value = 20835136
$ws.Range("L61").Value = 1868
$ws.Range("M61").Value = -20834924
$ws.Range("N61").Value = -2292

$ws.Range("H128").Value = 40430
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()

$ws.Range("H136").Value = 17545672
$ws.Range("I136").Value = 20835136
$ws.Range("J136").Value = 1868
$ws.Range("K136").Value = 62505408
$ws.Range("L136").Value = 5604
$ws.Range("M136").Value = -62502858
$ws.Range("N136").Value = -10704

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 18690.465
$ws.Range("I134").Value = 20769.72
$ws.Range("J134").Value = 1363.3334
$ws.Range("K134").Value = 62309.16
$ws.Range("L134").Value = 4090.0002
$ws.Range("M134").Value = -59774.16
$ws.Range("N134").Value = -9160.0002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 104.666664
$ws.Range("I7").Value = 88.2
$ws.Range("J7").Value = 125.25
$ws.Range("K7").Value = 88.2
$ws.Range("L7").Value = 125.25
$ws.Range("M7").Value = 24.8
$ws.Range("N7").Value = -351.25

$ws.Range("H31").Value = 1462.8387
$ws.Range("I31").Value = 1333.3636
$ws.Range("J31").Value = 1779.3334
$ws.Range("K31").Value = 1333.3636
$ws.Range("L31").Value = 1779.3334
$ws.Range("M31").Value = -1038.3636
$ws.Range("N31").Value = -2369.3334

$ws.Range("H34").Value = 1462.8387
$ws.Range("I34").Value = 1333.3636
$ws.Range("J34").Value = 1779.3334
$ws.Range("K34").Value = 1333.3636
$ws.Range("L34").Value = 1779.3334
$ws.Range("M34").Value = -1131.3636
$ws.Range("N34").Value = -2183.3334

$ws.Range("H132").Value = 3437.5652
$ws.Range("I132").Value = 3414.7778
$ws.Range("J132").Value = 3519.6
$ws.Range("K132").Value = 10244.3334
$ws.Range("L132").Value = 10558.8
$ws.Range("M132").Value = -7714.3334
$ws.Range("N132").Value = -15618.8

$ws.Range("H134").Value = 3548.5
$ws.Range("I134").Value = 3988.6
$ws.Range("J134").Value = 1348
$ws.Range("K134").Value = 11965.8
$ws.Range("L134").Value = 4044
$ws.Range("M134").Value = -9430.799999999999
$ws.Range("N134").Value = -9114

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 17545736
$ws.Range("I122").Value = 41667268
$ws.Range("J122").Value = 2804.4546
$ws.Range("K122").Value = 375005412
$ws.Range("L122").Value = 25240.0914
$ws.Range("M122").Value = -375002962
$ws.Range("N122").Value = -30140.0914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2198.8076
$ws.Range("I132").Value = 1851
$ws.Range("K132").Value = 5553
$ws.Range("M132").Value = -3023

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1457.1724
$ws.Range("I16").Value = 1629.76
$ws.Range("J16").Value = 378.5
$ws.Range("K16").Value = 1629.76
$ws.Range("L16").Value = 378.5
$ws.Range("M16").Value = -1459.76
$ws.Range("N16").Value = -718.5

$ws.Range("H46").Value = 1376.6666
$ws.Range("I46").Value = 1278
$ws.Range("J46").Value = 1500
$ws.Range("K46").Value = 1278
$ws.Range("L46").Value = 1500
$ws.Range("M46").Value = -1090
$ws.Range("N46").Value = -1876

$ws.Range("H61").Value = 1674.3
$ws.Range("I61").Value = 1487.6
$ws.Range("J61").Value = 1861
$ws.Range("K61").Value = 1487.6
$ws.Range("L61").Value = 1861
$ws.Range("M61").Value = -1285.6
$ws.Range("N61").Value = -2265

$ws.Range("H113").Value = 1674.3
$ws.Range("I113").Value = 1487.6
$ws.Range("J113").Value = 1861
$ws.Range("K113").Value = 1487.6
$ws.Range("L113").Value = 1861
$ws.Range("M113").Value = 682.4000000000001
$ws.Range("N113").Value = -6201

$ws.Range("H122").Value = 3539.48
$ws.Range("I122").Value = 3285.0952
$ws.Range("J122").Value = 4875
$ws.Range("K122").Value = 9855.285600000001
$ws.Range("L122").Value = 14625
$ws.Range("M122").Value = -7405.285600000001
$ws.Range("N122").Value = -19525

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4978.5713
$ws.Range("I81").Value = 2600
$ws.Range("J81").Value = 5375
$ws.Range("K81").Value = 5200
$ws.Range("L81").Value = 10750
$ws.Range("M81").Value = -4139
$ws.Range("N81").Value = -12872

$ws.Range("H84").Value = 4978.5713
$ws.Range("I84").Value = 2600
$ws.Range("J84").Value = 5375
$ws.Range("K84").Value = 26000
$ws.Range("L84").Value = 53750
$ws.Range("M84").Value = -20696
$ws.Range("N84").Value = -64358

$ws.Range("H136").Value = 20091.637
$ws.Range("I136").Value = 30429
$ws.Range("J136").Value = 2001.25
$ws.Range("K136").Value = 91287
$ws.Range("L136").Value = 6003.75
$ws.Range("M136").Value = -88737
$ws.Range("N136").Value = -11103.75
